$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the sheet's default column width to match the info file's refreshed
# metrics (matches the tiny defaultColWidth delta recorded for this sheet).
$ws.StandardWidth = 8.60546875

# Rename the sample (affects the shared strings used by rows 2-9 in columns A/B)
$ws.Range("A2:A5").Value = "cerebellum_tile1"
$ws.Range("B2:B5").Value = "039_GMB_tileRingMixScan_4rings_7scans"

# Clear out the now-unused data rows 6-9 (columns A,B,C,D,E) - only the D column
# keeps its cell (with its existing style) but with no value.
$ws.Range("A6:E9").ClearContents()

# Update the active selection
$ws.Range("A3:A5").Select()
